# New crime data collected — weekly CompStat report refresh (30th Precinct).
# Updates the report issue number / date range in the header, and refreshes
# the Crime Complaints table (rows 15-31) with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: volume/issue number and the "week covering" date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# ---------------------------------------------------------------------------
# Bulk numeric updates: cells whose style/number-format is unchanged, only
# the value moved.
# ---------------------------------------------------------------------------
$updates = @{
    "M15" = -60
    "N15" = -75.757575757575
    "C16" = 1
    "D16" = 3
    "E16" = -66.666666666666
    "F16" = 8
    "G16" = 10
    "H16" = -20
    "I16" = 57
    "J16" = 76
    "K16" = -25
    "L16" = -24
    "M16" = -55.46875
    "N16" = -86.198547215496
    "C17" = 1
    "E17" = -83.333333333333
    "F17" = 11
    "G17" = 20
    "H17" = -45
    "I17" = 95
    "J17" = 138
    "K17" = -31.159420289855
    "L17" = -31.654676258992
    "M17" = -15.929203539823
    "N17" = -78.841870824053
    "D18" = 1
    "E18" = 100
    "I18" = 44
    "J18" = 55
    "K18" = -20
    "L18" = 12.820512820512
    "M18" = -21.428571428571
    "N18" = -90.112359550561
    "C19" = 7
    "D19" = 9
    "E19" = -22.222222222222
    "F19" = 24
    "G19" = 22
    "H19" = 9.090909090909
    "I19" = 218
    "J19" = 209
    "K19" = 4.306220095693
    "L19" = 9.547738693467
    "M19" = 147.727272727273
    "N19" = 12.953367875647
    "C20" = 4
    "E20" = 300
    "F20" = 5
    "H20" = -28.571428571428
    "I20" = 45
    "J20" = 40
    "K20" = 12.5
    "L20" = -33.823529411764
    "M20" = 50
    "N20" = -75.274725274725
    "C21" = 15
    "D21" = 20
    "E21" = -25
    "F21" = 57
    "G21" = 67
    "H21" = -14.925373134328
    "I21" = 468
    "J21" = 527
    "K21" = -11.195445920303
    "L21" = -11.026615969581
    "M21" = 6.605922551252
    "N21" = -73.272415762421
    "F22" = 4
    "G22" = 2
    "H22" = 100
    "I22" = 14
    "K22" = 7.692307692307
    "L22" = -22.222222222222
    "M22" = 100
    "C24" = 10
    "D24" = 25
    "E24" = -60
    "F24" = 68
    "G24" = 71
    "H24" = -4.225352112676
    "I24" = 465
    "J24" = 453
    "K24" = 2.649006622516
    "L24" = -3.925619834710
    "M24" = 122.488038277512
    "C25" = 1
    "D25" = 10
    "E25" = -90
    "G25" = 25
    "H25" = -68
    "I25" = 90
    "J25" = 123
    "K25" = -26.829268292682
    "L25" = -29.6875
    "C26" = 8
    "D26" = 3
    "E26" = 166.666666666667
    "G26" = 20
    "H26" = 65
    "I26" = 209
    "J26" = 196
    "K26" = 6.632653061224
    "L26" = 12.972972972973
    "M26" = -28.668941979522
    "I28" = 20
    "K28" = 42.857142857142
    "N29" = -96.491228070175
    "N30" = -96.363636363636
    "J31" = 2
    "K31" = -50
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

# ---------------------------------------------------------------------------
# Cells that flip from the blank/not-applicable text style ("0" / "***.*")
# to a real numeric entry now that data exists for them. Pull the number
# format from a sibling cell in the same column family, then set the value.
# ---------------------------------------------------------------------------
$ws.Range("C16").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1

$ws.Range("C16").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1

$ws.Range("E16").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100

$ws.Range("C16").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1

$ws.Range("E16").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100

# ---------------------------------------------------------------------------
# Row 22 (Transit): the Week-to-Date 2024 count (D22) and %Chg (E22) flip
# the other way — from real numbers to the "no prior data" placeholder text.
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"

$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$excel.CutCopyMode = $false
